# The interactions dataframes were missing a few derived columns
# (TD, GP_N, GP_T, GP_T_SC) that need to sit right after SCENARIO and
# before the existing OI_* columns. Insert 4 new columns at F:I, which
# shifts the old F:S block to J:W, then populate the new columns and
# refresh the OI_T / OI_T_SS / OI_T_GR values (now U/V/W) that changed
# once the pairing bug was fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank columns before column F (old F:I -> new J:M, etc.)
$ws.Range("F1:I1").EntireColumn.Insert()

# New header row for the inserted columns
$ws.Range("F1").Value = "TD"
$ws.Range("G1").Value = "GP_N"
$ws.Range("H1").Value = "GP_T"
$ws.Range("I1").Value = "GP_T_SC"

# New data for rows 2-6
$ws.Range("F2").Value = 1433.02
$ws.Range("G2").Value = 24
$ws.Range("H2").Value = 57.91
$ws.Range("I2").Value = 49.55

$ws.Range("F3").Value = 282.26
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 47.28
$ws.Range("I3").Value = 45.36

$ws.Range("F4").Value = 199.31
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 42.27
$ws.Range("I4").Value = 35.11

$ws.Range("F5").Value = 385.14
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 68.02
$ws.Range("I5").Value = 68.18000000000001

$ws.Range("F6").Value = 545.02
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 43.86
$ws.Range("I6").Value = "NULL"

# Corrected OI_T / OI_T_SS / OI_T_GR values (now columns U/V/W after the shift)
$ws.Range("U2").Value = 8.74
$ws.Range("V2").Value = 0.45
$ws.Range("W2").Value = 20.79

$ws.Range("U3").Value = 3.03
$ws.Range("V3").Value = 0.43
$ws.Range("W3").Value = 2.07

$ws.Range("U4").Value = 1.52
$ws.Range("V4").Value = 0.16
$ws.Range("W4").Value = 4.6

$ws.Range("U5").Value = 6.8
$ws.Range("V5").Value = 0.5
$ws.Range("W5").Value = 21.21

$ws.Range("U6").Value = 10.36
$ws.Range("V6").Value = 0.37
$ws.Range("W6").Value = 47.96
